$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.876.68"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "3.139.55"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "3.132.64"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("E13").Value = "  -2.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "3.659.22"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").Value = "63.740.85"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "3.137.79"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("E26").Value = "  +6.72%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.83%  "
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("D35").Value = "0.0₃0846"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.28%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "453.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  +4.85%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "2.913.49"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.61%  "
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.69%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.111"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.79%  "
